$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Армагедон"
$ws.Range("B13").Value = "entry.934056410"

$ws.Range("A14").Value = "Електроопора"
$ws.Range("B14").Value = "entry.1240912437"

$ws.Range("A15").Value = "Комплімент"
$ws.Range("B15").Value = "entry.2030694513"

$ws.Range("B13").Select()
